$wb = $excel.ActiveWorkbook

# ALC!row32
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 29816.75
$ws.Range("I32").Value = 211499.5
$ws.Range("J32").Value = 3862.0715
$ws.Range("K32").Value = 211499.5
$ws.Range("L32").Value = 3862.0715
$ws.Range("M32").Value = -211173.5
$ws.Range("N32").Value = -4514.0715

# ALC!row54
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("N54").ClearContents()

# ALC!row97
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 737.4
$ws.Range("J97").Value = 771.75
$ws.Range("L97").Value = 2315.25
$ws.Range("N97").Value = -3307.25

# ALC!row111
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 22290.637
$ws.Range("I111").Value = 1689.8334
$ws.Range("J111").Value = 47011.6
$ws.Range("K111").Value = 5069.5002
$ws.Range("L111").Value = 141034.8
$ws.Range("M111").Value = -2002.5002
$ws.Range("N111").Value = -147168.8

# ALC!row116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 36530296
$ws.Range("I116").Value = 50202280
$ws.Range("K116").Value = 50202280
$ws.Range("M116").Value = -50198838

# ALC!row131
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H131").Value = 14640
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 14640
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 43920
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -54000

# ALC!row132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4713.567
$ws.Range("I132").Value = 4265.6924
$ws.Range("J132").Value = 7624.75
$ws.Range("K132").Value = 12797.0772
$ws.Range("L132").Value = 22874.25
$ws.Range("M132").Value = -10267.0772
$ws.Range("N132").Value = -27934.25

# ALC!row137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1440.4242
$ws.Range("I137").Value = 1426.862
$ws.Range("J137").Value = 1538.75
$ws.Range("K137").Value = 4280.586
$ws.Range("L137").Value = 4616.25
$ws.Range("M137").Value = -1730.586
$ws.Range("N137").Value = -9716.25

# ARM!row2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1389.2727
$ws.Range("I2").Value = 882.6667
$ws.Range("K2").Value = 882.6667
$ws.Range("M2").Value = -769.6667

# ARM!row116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1389.2727
$ws.Range("I116").Value = 882.6667
$ws.Range("K116").Value = 882.6667
$ws.Range("M116").Value = 1411.3333

# BSM!row3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1389.2727
$ws.Range("I3").Value = 882.6667
$ws.Range("K3").Value = 882.6667
$ws.Range("M3").Value = -768.6667

# BSM!row105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2040.5555
$ws.Range("I105").Value = 1966.5
$ws.Range("J105").Value = 2188.6667
$ws.Range("K105").Value = 1966.5
$ws.Range("L105").Value = 2188.6667
$ws.Range("M105").Value = -219.5
$ws.Range("N105").Value = -5682.6667

# BSM!row107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 27802994
$ws.Range("I107").Value = 16792.8
$ws.Range("J107").Value = 166734000
$ws.Range("K107").Value = 16792.8
$ws.Range("L107").Value = 166734000
$ws.Range("M107").Value = -14872.8
$ws.Range("N107").Value = -166737840

# CRP!row16
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4999.6665

# CRP!row31
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1862.125
$ws.Range("I31").Value = 1079.3
$ws.Range("J31").Value = 3166.8333
$ws.Range("K31").Value = 1079.3
$ws.Range("L31").Value = 3166.8333
$ws.Range("M31").Value = -784.3
$ws.Range("N31").Value = -3756.8333

# CRP!row34
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1862.125
$ws.Range("I34").Value = 1079.3
$ws.Range("J34").Value = 3166.8333
$ws.Range("K34").Value = 1079.3
$ws.Range("L34").Value = 3166.8333
$ws.Range("M34").Value = -877.3
$ws.Range("N34").Value = -3570.8333

# CRP!row99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3462.8572
$ws.Range("I99").Value = 3413
$ws.Range("K99").Value = 3413
$ws.Range("M99").Value = -1915

# CRP!row107
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1528.1578
$ws.Range("I107").Value = 1474.3214
$ws.Range("J107").Value = 1678.9
$ws.Range("K107").Value = 1474.3214
$ws.Range("L107").Value = 1678.9
$ws.Range("M107").Value = 445.6786
$ws.Range("N107").Value = -5518.9

# CRP!row113
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 4999.6665

# CRP!row126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 3462.8572
$ws.Range("I126").Value = 3413
$ws.Range("K126").Value = 10239
$ws.Range("M126").Value = -7769

# CRP!row132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2230.1052
$ws.Range("J132").Value = 5885.5
$ws.Range("L132").Value = 17656.5
$ws.Range("N132").Value = -22716.5

# CRP!row134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1580.85
$ws.Range("I134").Value = 1223.2222
$ws.Range("K134").Value = 3669.6666
$ws.Range("M134").Value = -1134.6666

# CUL!row9
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 4569.857
$ws.Range("I9").Value = 4397.8
$ws.Range("J9").Value = 5000
$ws.Range("K9").Value = 13193.4
$ws.Range("L9").Value = 15000
$ws.Range("M9").Value = -12969.4
$ws.Range("N9").Value = -15448

# CUL!row13
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 488.63635
$ws.Range("I13").Value = 99
$ws.Range("J13").Value = 813.3333
$ws.Range("K13").Value = 297
$ws.Range("L13").Value = 2439.9999
$ws.Range("M13").Value = -129
$ws.Range("N13").Value = -2775.9999

# CUL!row57
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 5498.5
$ws.Range("I57").Value = 3999
$ws.Range("J57").Value = 6998
$ws.Range("K57").Value = 11997
$ws.Range("L57").Value = 20994
$ws.Range("M57").Value = -11438
$ws.Range("N57").Value = -22112

# CUL!row62
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 5125
$ws.Range("I62").Value = 2000
$ws.Range("J62").Value = 6166.6665
$ws.Range("K62").Value = 6000
$ws.Range("L62").Value = 18499.9995
$ws.Range("M62").Value = -5314
$ws.Range("N62").Value = -19871.9995

# CUL!row65
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H65").Value = 5125
$ws.Range("I65").Value = 2000
$ws.Range("J65").Value = 6166.6665
$ws.Range("K65").Value = 18000
$ws.Range("L65").Value = 55499.9985
$ws.Range("M65").Value = -14568
$ws.Range("N65").Value = -62363.9985

# CUL!row112
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 2582.8333
$ws.Range("I112").Value = 2919.4
$ws.Range("K112").Value = 8758.200000000001
$ws.Range("M112").Value = -7650.200000000001

# GSM!row70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9857.208000000001
$ws.Range("J70").Value = 12419.571
$ws.Range("L70").Value = 12419.571
$ws.Range("N70").Value = -12959.571

# GSM!row73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 9857.208000000001
$ws.Range("J73").Value = 12419.571
$ws.Range("L73").Value = 12419.571
$ws.Range("N73").Value = -14291.571

# GSM!row80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3624.2666
$ws.Range("I80").Value = 4082.375
$ws.Range("J80").Value = 3100.7144
$ws.Range("K80").Value = 4082.375
$ws.Range("L80").Value = 3100.7144
$ws.Range("M80").Value = -3084.375
$ws.Range("N80").Value = -5096.7144

# GSM!row83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 3624.2666
$ws.Range("I83").Value = 4082.375
$ws.Range("J83").Value = 3100.7144
$ws.Range("K83").Value = 20411.875
$ws.Range("L83").Value = 15503.572
$ws.Range("M83").Value = -15419.875
$ws.Range("N83").Value = -25487.572

# GSM!row102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1809.7778
$ws.Range("I102").Value = 1448.2858
$ws.Range("J102").Value = 3075
$ws.Range("K102").Value = 1448.2858
$ws.Range("L102").Value = 3075
$ws.Range("M102").Value = 173.7141999999999
$ws.Range("N102").Value = -6319

# GSM!row113
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1866.75
$ws.Range("I113").Value = 1326.8
$ws.Range("J113").Value = 2766.6667
$ws.Range("K113").Value = 1326.8
$ws.Range("L113").Value = 2766.6667
$ws.Range("M113").Value = 843.2
$ws.Range("N113").Value = -7106.6667

# GSM!row126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 6805.615
$ws.Range("I126").Value = 7687.2856
$ws.Range("K126").Value = 23061.8568
$ws.Range("M126").Value = -20591.8568

# LTW!row7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 62500824
$ws.Range("I7").Value = 71429384
$ws.Range("K7").Value = 71429384
$ws.Range("M7").Value = -71429272

# LTW!row46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1635.4688
$ws.Range("I46").Value = 783
$ws.Range("J46").Value = 1919.625
$ws.Range("K46").Value = 783
$ws.Range("L46").Value = 1919.625
$ws.Range("M46").Value = -595
$ws.Range("N46").Value = -2295.625

# LTW!row61
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 18551.715
$ws.Range("I61").Value = 16195.556
$ws.Range("J61").Value = 22792.8
$ws.Range("K61").Value = 16195.556
$ws.Range("L61").Value = 22792.8
$ws.Range("M61").Value = -15993.556
$ws.Range("N61").Value = -23196.8

# LTW!row113
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 18551.715
$ws.Range("I113").Value = 16195.556
$ws.Range("J113").Value = 22792.8
$ws.Range("K113").Value = 16195.556
$ws.Range("L113").Value = 22792.8
$ws.Range("M113").Value = -14025.556
$ws.Range("N113").Value = -27132.8

# LTW!row122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3084.4827
$ws.Range("I122").Value = 2212.1428
$ws.Range("J122").Value = 5374.375
$ws.Range("K122").Value = 6636.428400000001
$ws.Range("L122").Value = 16123.125
$ws.Range("M122").Value = -4186.428400000001
$ws.Range("N122").Value = -21023.125

# LTW!row126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 62500824
$ws.Range("I126").Value = 71429384
$ws.Range("K126").Value = 214288152
$ws.Range("M126").Value = -214285682

# LTW!row132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4477.879
$ws.Range("I132").Value = 2398.375
$ws.Range("J132").Value = 10023.223
$ws.Range("K132").Value = 7195.125
$ws.Range("L132").Value = 30069.669
$ws.Range("M132").Value = -4665.125
$ws.Range("N132").Value = -35129.669

# LTW!row136
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2045.7576
$ws.Range("I136").Value = 1785.1052
$ws.Range("K136").Value = 5355.3156
$ws.Range("M136").Value = -2805.3156

# WVR!row107
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 571.25
$ws.Range("J107").Value = 598.5
$ws.Range("L107").Value = 1795.5
$ws.Range("N107").Value = -5635.5

# WVR!row122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2139.5
$ws.Range("I122").Value = 2256.7334
$ws.Range("K122").Value = 6770.2002
$ws.Range("M122").Value = -4320.2002

# WVR!row126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1352.16
$ws.Range("I126").Value = 1308.5
$ws.Range("K126").Value = 3925.5
$ws.Range("M126").Value = -1455.5

# WVR!row132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6995.7144
$ws.Range("I132").Value = 7867.625
$ws.Range("K132").Value = 23602.875
$ws.Range("M132").Value = -21072.875
